$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header block (rows 5-9) ---
# ATTN value cleared (was "PT. KLINE TOTAL LOGISTICS INDONESIA")
$ws.Range("C5").Value = ""
# FROM value updated
$ws.Range("C6").Value = "WISNU - FDK INDONESIA"
# PPBE NO. value updated
$ws.Range("J7").Value = "130/W/20"

# --- Data block ---
# Row 12: first container line item
$ws.Range("A12").Value = "FIENR20-015"
$ws.Range("B12").Value = "EVEREADY GOLD US LR6 BULK (E7075601)"
$ws.Range("C12").Value = 52726
$ws.Range("D12").Value = "E7075601"
$ws.Range("E12").Value = 753920
$ws.Range("F12").Value = 1216
$ws.Range("G12").Value = 18452.8
$ws.Range("H12").Value = 17341.3
$ws.Range("I12").Value = 14.1792
$ws.Range("J12").Value = 19
$ws.Range("K12").Value = "TGHU 5244 291 / ZZC-SB 109 121" + [char]10 + "(40 FEET)"
$ws.Range("L12").Value = 3700
$ws.Range("M12").Value = 22152.8

# Row 13: new sub-total row inserted between the two container lines.
# Clone formatting from the existing sub-total row (row 15) first, then
# overwrite with the correct content.
$ws.Range("A15:M15").Copy()
$ws.Range("A13:M13").PasteSpecial(-4122)
$ws.Range("A13").Value = ""
$ws.Range("B13").Value = ""
$ws.Range("C13").Value = "TOTAL"
$ws.Range("D13").Value = ""
$ws.Range("E13").Value = 753920
$ws.Range("F13").Value = 1216
$ws.Range("G13").Value = 18452.8
$ws.Range("H13").Value = 17341.3
$ws.Range("I13").Value = 14.179
$ws.Range("J13").Value = 19
$ws.Range("K13").Value = ""
$ws.Range("L13").Value = 3700
$ws.Range("M13").Value = 22152.8
$ws.Range("A13:B13").Merge()
$ws.Range("C13:D13").Merge()

# Row 14: second container line item
$ws.Range("A14").Value = "FIENR20-015"
$ws.Range("B14").Value = "EVEREADY GOLD US LR6 BULK (E7075601)"
$ws.Range("C14").Value = 52726
$ws.Range("D14").Value = "E7075601"
$ws.Range("E14").Value = 753920
$ws.Range("F14").Value = 1216
$ws.Range("G14").Value = 18452.8
$ws.Range("H14").Value = 17341.3
$ws.Range("I14").Value = 14.1792
$ws.Range("J14").Value = 19
$ws.Range("K14").Value = "TGHU 5245 488 / ZZC-SB 111 849" + [char]10 + "(40 FEET)"
$ws.Range("L14").Value = 3700
$ws.Range("M14").Value = 22152.8

# Row 15: existing sub-total row, refreshed totals (label "TOTAL" unchanged)
$ws.Range("E15").Value = 753920
$ws.Range("F15").Value = 1216
$ws.Range("G15").Value = 18452.8
$ws.Range("H15").Value = 17341.3
$ws.Range("I15").Value = 14.179
$ws.Range("J15").Value = 19
$ws.Range("L15").Value = 3700
$ws.Range("M15").Value = 22152.8

# Row 17: grand-total row, refreshed totals (label "TOTAL ALL" unchanged)
$ws.Range("E17").Value = 1507840
$ws.Range("F17").Value = 2432
$ws.Range("G17").Value = 36905.6
$ws.Range("H17").Value = 34682.6
$ws.Range("I17").Value = 28.358
$ws.Range("J17").Value = 38
$ws.Range("L17").Value = 7400
$ws.Range("M17").Value = 44305.6

# --- Column widths (best-fit refresh following the wider new text) ---
$ws.Columns.Item(1).ColumnWidth = 13.996582 - 0.8333333333333334
$ws.Columns.Item(2).ColumnWidth = 43.560791 - 0.8333333333333334
$ws.Columns.Item(5).ColumnWidth = 16.567383 - 0.8333333333333334
$ws.Columns.Item(7).ColumnWidth = 12.568359 - 0.8333333333333334
$ws.Columns.Item(8).ColumnWidth = 12.568359 - 0.8333333333333334
$ws.Columns.Item(11).ColumnWidth = 36.419678 - 0.8333333333333334
$ws.Columns.Item(13).ColumnWidth = 12.568359 - 0.8333333333333334
